# Weekly driver report update for 2025-04-28
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Bad Drivers" summary row (row 3) and Totals row (row 4) ---
$ws.Range("B3").Value2 = 4
$ws.Range("C3").Value2 = 214
$ws.Range("D3").Value2 = 97.3

$ws.Range("B4").Value2 = 4
$ws.Range("C4").Value2 = 214

# --- "Good Drivers" table: a brand-new driver (21.40.1.3) showed up this
#     week, so insert a fresh row at the top of the table (row 12) and let
#     the rest of the rows shift down by one. ---
$ws.Rows(12).Insert()

# The inserted row doesn't automatically pick up the table's row styling,
# so copy formats down from the row right below it (now row 13, which used
# to be row 12) before filling in the new driver's numbers.
$ws.Range("A13:E13").Copy()
$ws.Range("A12").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Range("A12").Value2 = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.40.1.3"
$ws.Range("B12").Value2 = 11128
$ws.Range("D12").Value2 = 100
$ws.Range("E12").Value2 = 0

# Refresh the sample counts for the rest of the table now that the weekly
# numbers have rolled forward (driver versions/dates stay put - they just
# shifted down a row with the insert above).
$ws.Range("B13").Value2 = 486214

$ws.Range("B14").Value2 = 79953

$ws.Range("B15").Value2 = 35355

$ws.Range("B16").Value2 = 65425

$ws.Range("B17").Value2 = 117653

Write-Host "done"
